$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.516.31"
$ws.Range("E2").Value = "'  +2.12%  "
$ws.Range("D3").Value = "'1.680.84"
$ws.Range("E3").Value = "'  +3.02%  "
$ws.Range("E4").Value = "'  +0.09%  "
$ws.Range("D5").Value = "'216.66"
$ws.Range("E5").Value = "'  +2.76%  "
$ws.Range("D6").Value = "'0.5324"
$ws.Range("E6").Value = "'  +1.82%  "
$ws.Range("E7").Value = "'  +0.10%  "
$ws.Range("D8").Value = "'0.2683"
$ws.Range("E8").Value = "'  +4.01%  "
$ws.Range("D9").Value = "'0.06399"
$ws.Range("E9").Value = "'  +2.01%  "
$ws.Range("E10").Value = "'  +5.43%  "
$ws.Range("D11").Value = "'0.07798"
$ws.Range("E11").Value = "'  +2.85%  "
$ws.Range("D12").Value = "'1.680.13"
$ws.Range("E12").Value = "'  +3.04%  "
$ws.Range("D13").Value = "'4.499"
$ws.Range("E13").Value = "'  +2.00%  "
$ws.Range("E14").Value = "'  +1.43%  "
$ws.Range("D15").Value = "'0.0₅8331"
$ws.Range("E15").Value = "'  +3.99%  "
$ws.Range("D16").Value = "'65.62"
$ws.Range("E16").Value = "'  +1.48%  "
$ws.Range("D17").Value = "'26.551.62"
$ws.Range("E17").Value = "'  +2.34%  "
$ws.Range("E18").Value = "'  -0.03%  "
$ws.Range("D19").Value = "'4.757"
$ws.Range("E19").Value = "'  +1.87%  "
$ws.Range("D20").Value = "'194.49"
$ws.Range("E20").Value = "'  +4.92%  "
$ws.Range("E21").Value = "'  +2.28%  "
$ws.Range("D22").Value = "'6.348"
$ws.Range("E22").Value = "'  +3.93%  "
$ws.Range("E23").Value = "'  +0.10%  "
$ws.Range("D24").Value = "'143.39"
$ws.Range("E24").Value = "'  -1.29%  "
$ws.Range("D25").Value = "'0.1283"
$ws.Range("E25").Value = "'  +5.81%  "
$ws.Range("D26").Value = "'7.432"
$ws.Range("E26").Value = "'  +0.71%  "
$ws.Range("D27").Value = "'16.32"
$ws.Range("E27").Value = "'  +4.20%  "
$ws.Range("D28").Value = "'1.427"
$ws.Range("E28").Value = "'  +4.18%  "
$ws.Range("D29").Value = "'0.06217"
$ws.Range("E29").Value = "'  +5.45%  "
$ws.Range("E30").Value = "'  +2.47%  "
$ws.Range("E31").Value = "'  +5.41%  "
$ws.Range("D32").Value = "'3.455"
$ws.Range("E32").Value = "'  +1.98%  "
$ws.Range("D33").Value = "'1.690"
$ws.Range("E33").Value = "'  +4.13%  "
$ws.Range("E34").Value = "'  +2.82%  "
$ws.Range("D35").Value = "'2.425"
$ws.Range("E35").Value = "'  +1.69%  "
$ws.Range("D36").Value = "'2.789"
$ws.Range("E36").Value = "'  +1.27%  "
$ws.Range("D37").Value = "'0.5726"
$ws.Range("E37").Value = "'  -0.94%  "
$ws.Range("D38").Value = "'0.01637"
$ws.Range("E38").Value = "'  +2.29%  "
$ws.Range("D39").Value = "'6.032"
$ws.Range("E39").Value = "'  +6.19%  "
$ws.Range("D40").Value = "'1.075.20"
$ws.Range("E40").Value = "'  +3.81%  "
$ws.Range("D41").Value = "'0.8602"
$ws.Range("E41").Value = "'  +1.69%  "
$ws.Range("E42").Value = "'  -0.32%  "
$ws.Range("D43").Value = "'99.97"
$ws.Range("E43").Value = "'  -0.03%  "
$ws.Range("D44").Value = "'1.827.05"
$ws.Range("D45").Value = "'0.0₈108"
$ws.Range("E45").Value = "'  -0.76%  "
$ws.Range("D46").Value = "'57.06"
$ws.Range("E46").Value = "'  +3.92%  "
$ws.Range("D47").Value = "'8.138"
$ws.Range("E47").Value = "'  +1.23%  "
$ws.Range("E48").Value = "'  +0.26%  "
$ws.Range("D49").Value = "'0.05207"
$ws.Range("E49").Value = "'  +0.86%  "
$ws.Range("E50").Value = "'  +5.99%  "
$ws.Range("D51").Value = "'6.024"
$ws.Range("E51").Value = "'  +2.95%  "
